$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 27778014
$ws.Range("I6").Value = 33333462
$ws.Range("K6").Value = 100000386
$ws.Range("M6").Value = -100000274
$ws.Range("H113").Value = 31342.084
$ws.Range("I113").Value = 34888
$ws.Range("J113").Value = 24250.25
$ws.Range("K113").Value = 34888
$ws.Range("L113").Value = 24250.25
$ws.Range("M113").Value = -31634
$ws.Range("N113").Value = -30758.25
$ws.Range("H137").Value = 454760.7
$ws.Range("I137").Value = 738223.4
$ws.Range("J137").Value = 13818.777
$ws.Range("K137").Value = 2214670.2
$ws.Range("L137").Value = 41456.331
$ws.Range("M137").Value = -2212120.2
$ws.Range("N137").Value = -46556.331

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 168444.92
$ws.Range("I45").Value = 240610.11
$ws.Range("J45").Value = 6073.25
$ws.Range("K45").Value = 240610.11
$ws.Range("L45").Value = 6073.25
$ws.Range("M45").Value = -240233.11
$ws.Range("N45").Value = -6827.25
$ws.Range("H61").Value = 8463.571
$ws.Range("I61").Value = 8749.166999999999
$ws.Range("J61").Value = 6750
$ws.Range("K61").Value = 8749.166999999999
$ws.Range("L61").Value = 6750
$ws.Range("M61").Value = -8537.166999999999
$ws.Range("N61").Value = -7174
$ws.Range("H74").Value = 2388.1785
$ws.Range("I74").Value = 2206.5
$ws.Range("J74").Value = 4750
$ws.Range("K74").Value = 2206.5
$ws.Range("L74").Value = 4750
$ws.Range("M74").Value = -1332.5
$ws.Range("N74").Value = -6498
$ws.Range("H77").Value = 2388.1785
$ws.Range("I77").Value = 2206.5
$ws.Range("J77").Value = 4750
$ws.Range("K77").Value = 11032.5
$ws.Range("L77").Value = 23750
$ws.Range("M77").Value = -6664.5
$ws.Range("N77").Value = -32486
$ws.Range("H102").Value = 5516.553
$ws.Range("I102").Value = 3934.2666
$ws.Range("J102").Value = 8308.823
$ws.Range("K102").Value = 3934.2666
$ws.Range("L102").Value = 8308.823
$ws.Range("M102").Value = -2312.2666
$ws.Range("N102").Value = -11552.823
$ws.Range("H132").Value = 2384.484
$ws.Range("I132").Value = 1859.2759
$ws.Range("K132").Value = 5577.8277
$ws.Range("M132").Value = -3047.8277
$ws.Range("H136").Value = 8463.571
$ws.Range("I136").Value = 8749.166999999999
$ws.Range("J136").Value = 6750
$ws.Range("K136").Value = 26247.501
$ws.Range("L136").Value = 20250
$ws.Range("M136").Value = -23697.501
$ws.Range("N136").Value = -25350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7002.885
$ws.Range("I134").Value = 8004
$ws.Range("J134").Value = 4285.5713
$ws.Range("K134").Value = 24012
$ws.Range("L134").Value = 12856.7139
$ws.Range("M134").Value = -21477
$ws.Range("N134").Value = -17926.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1926.5834
$ws.Range("I31").Value = 1310.9166
$ws.Range("J31").Value = 2337.0278
$ws.Range("K31").Value = 1310.9166
$ws.Range("L31").Value = 2337.0278
$ws.Range("M31").Value = -1015.9166
$ws.Range("N31").Value = -2927.0278
$ws.Range("H34").Value = 1926.5834
$ws.Range("I34").Value = 1310.9166
$ws.Range("J34").Value = 2337.0278
$ws.Range("K34").Value = 1310.9166
$ws.Range("L34").Value = 2337.0278
$ws.Range("M34").Value = -1108.9166
$ws.Range("N34").Value = -2741.0278
$ws.Range("H58").Value = 4211.173
$ws.Range("I58").Value = 3918.2092
$ws.Range("K58").Value = 3918.2092
$ws.Range("M58").Value = -3715.2092
$ws.Range("H122").Value = 1232.9333
$ws.Range("I122").Value = 972.5
$ws.Range("J122").Value = 1753.8
$ws.Range("K122").Value = 2917.5
$ws.Range("L122").Value = 5261.4
$ws.Range("M122").Value = -467.5
$ws.Range("N122").Value = -10161.4
$ws.Range("H136").Value = 4211.173
$ws.Range("I136").Value = 3918.2092
$ws.Range("K136").Value = 11754.6276
$ws.Range("M136").Value = -9204.6276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1740
$ws.Range("I17").Value = 1296.6666
$ws.Range("J17").Value = 3070
$ws.Range("K17").Value = 3889.9998
$ws.Range("L17").Value = 9210
$ws.Range("M17").Value = -3720.9998
$ws.Range("N17").Value = -9548
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 1500
$ws.Range("N32").Value = -2066
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 2274842.5
$ws.Range("J34").Value = 4133
$ws.Range("L34").Value = 12399
$ws.Range("N34").Value = -12567
$ws.Range("H39").Value = 3925
$ws.Range("J39").Value = 2437.5
$ws.Range("L39").Value = 7312.5
$ws.Range("N39").Value = -7900.5
$ws.Range("H52").Value = 57500
$ws.Range("J52").Value = 57500
$ws.Range("L52").Value = 172500
$ws.Range("N52").Value = -173032
$ws.Range("H55").Value = 5920.5
$ws.Range("I55").Value = 360.5
$ws.Range("J55").Value = 7310.5
$ws.Range("K55").Value = 1081.5
$ws.Range("L55").Value = 21931.5
$ws.Range("M55").Value = -904.5
$ws.Range("N55").Value = -22285.5
$ws.Range("H68").Value = 11634699
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 13896419
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 41689257
$ws.Range("M68").Value = -8189
$ws.Range("N68").Value = -41690879
$ws.Range("H71").Value = 11634699
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 13896419
$ws.Range("K71").Value = 27000
$ws.Range("L71").Value = 125067771
$ws.Range("M71").Value = -22944
$ws.Range("N71").Value = -125075883
$ws.Range("H81").Value = 10190.706
$ws.Range("I81").Value = 3704.3333
$ws.Range("J81").Value = 11580.643
$ws.Range("K81").Value = 11112.9999
$ws.Range("L81").Value = 34741.929
$ws.Range("M81").Value = -9989.999899999999
$ws.Range("N81").Value = -36987.929
$ws.Range("H84").Value = 10190.706
$ws.Range("I84").Value = 3704.3333
$ws.Range("J84").Value = 11580.643
$ws.Range("K84").Value = 33338.9997
$ws.Range("L84").Value = 104225.787
$ws.Range("M84").Value = -27722.9997
$ws.Range("N84").Value = -115457.787
$ws.Range("H107").Value = 5370.778
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 5917.125
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 17751.375
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -21591.375
$ws.Range("H115").Value = 3367.6667
$ws.Range("I115").Value = 444.25
$ws.Range("K115").Value = 1332.75
$ws.Range("M115").Value = -157.75
$ws.Range("H124").Value = 4180.077
$ws.Range("I124").Value = 923.625
$ws.Range("J124").Value = 9390.4
$ws.Range("K124").Value = 2770.875
$ws.Range("L124").Value = 28171.2
$ws.Range("M124").Value = 2139.125
$ws.Range("N124").Value = -37991.2
$ws.Range("H133").Value = 10811
$ws.Range("I133").Value = 10013.75
$ws.Range("K133").Value = 30041.25
$ws.Range("M133").Value = -24981.25
$ws.Range("H137").Value = 8176.5557
$ws.Range("I137").Value = 2890
$ws.Range("J137").Value = 12405.8
$ws.Range("K137").Value = 8670
$ws.Range("L137").Value = 37217.39999999999
$ws.Range("M137").Value = -3570
$ws.Range("N137").Value = -47417.39999999999
$ws.Range("H139").Value = 1306871.9
$ws.Range("I139").Value = 1430859.6
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 4292578.800000001
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -4287438.800000001
$ws.Range("N139").Value = -25280
$ws.Range("H141").Value = 1749.25
$ws.Range("I141").Value = 999
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 2997
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 2183
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 321514340
$ws.Range("J98").Value = 321514340
$ws.Range("L98").Value = 321514340
$ws.Range("N98").Value = -321520330
$ws.Range("H102").Value = 27596.268
$ws.Range("I102").Value = 25828.666
$ws.Range("J102").Value = 34666.668
$ws.Range("K102").Value = 25828.666
$ws.Range("L102").Value = 34666.668
$ws.Range("M102").Value = -24206.666
$ws.Range("N102").Value = -37910.668
$ws.Range("H132").Value = 2363.32
$ws.Range("I132").Value = 2079.8044
$ws.Range("J132").Value = 5623.75
$ws.Range("K132").Value = 6239.4132
$ws.Range("L132").Value = 16871.25
$ws.Range("M132").Value = -3709.4132
$ws.Range("N132").Value = -21931.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4007.842
$ws.Range("I122").Value = 2689.3333
$ws.Range("J122").Value = 7244.1816
$ws.Range("K122").Value = 8067.999899999999
$ws.Range("L122").Value = 21732.5448
$ws.Range("M122").Value = -5617.999899999999
$ws.Range("N122").Value = -26632.5448
$ws.Range("H132").Value = 30735.143
$ws.Range("I132").Value = 37042.867
$ws.Range("J132").Value = 14965.833
$ws.Range("K132").Value = 111128.601
$ws.Range("L132").Value = 44897.499
$ws.Range("M132").Value = -108598.601
$ws.Range("N132").Value = -49957.499
$ws.Range("H136").Value = 4598.5713
$ws.Range("I136").Value = 4407.1665
$ws.Range("J136").Value = 5747
$ws.Range("K136").Value = 13221.4995
$ws.Range("L136").Value = 17241
$ws.Range("M136").Value = -10671.4995
$ws.Range("N136").Value = -22341
